# "Generate Report for Archive"
# Status text moved on from handoff -> the items are now actively being
# translated, so every "Ready for handoff" status cell becomes
# "In Translation" across the Overview summary sheet and each of the two
# per-locale detail sheets. The Status columns are then re-sized to fit
# the (shorter) new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: per-locale status columns E (zh-cn) and F (de-de) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

# --- zh-cn detail sheet: Status column C ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

# --- de-de detail sheet: Status column C ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# Re-fit the status columns now that the text is shorter than
# "Ready for handoff" was.
$overview.Columns.Item(5).ColumnWidth = 13.4101845877511
$overview.Columns.Item(6).ColumnWidth = 13.4101845877511
$zhcn.Columns.Item(3).ColumnWidth = 13.4101845877511
$dede.Columns.Item(3).ColumnWidth = 13.4101845877511
